$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# Cells in columns D (Price) and E (Volume/1h) hold numeric-looking text that must
# stay as literal text (matches original inlineStr cells, preserves exact digits/
# trailing zeros and the "%" formatting). Force the Text number format first so
# Excel does not silently convert the assigned string into a float.
$textCells = @(
    "D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "D19", "E19", "D20", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "E25", "D38", "E38", "D39", "E39", "D40", "E40", "D41", "E41", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "D47", "E47", "D48", "E48", "D49", "E49", "D50", "E50", "D51", "E51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "297.33"
$ws.Range("E2").Value = "1.75%"
$ws.Range("D3").Value = "41.95"
$ws.Range("E3").Value = "4.00%"
$ws.Range("D4").Value = "5.011"
$ws.Range("E4").Value = "-0.36%"
$ws.Range("D5").Value = "0.07522"
$ws.Range("E5").Value = "2.70%"
$ws.Range("D6").Value = "1.583"
$ws.Range("E6").Value = "3.23%"
$ws.Range("D7").Value = "0.9259"
$ws.Range("E7").Value = "-0.16%"
$ws.Range("D9").Value = "0.1195"
$ws.Range("E9").Value = "0.65%"
$ws.Range("D10").Value = "0.1828"
$ws.Range("E10").Value = "4.79%"
$ws.Range("D11").Value = "0.08935"
$ws.Range("E11").Value = "3.25%"
$ws.Range("D12").Value = "0.04079"
$ws.Range("E12").Value = "-5.75%"
$ws.Range("E13").Value = "-0.53%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "0.001280"
$ws.Range("E14").Value = "0.41%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "0.005838"
$ws.Range("E15").Value = "-2.78%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "3.356"
$ws.Range("E16").Value = "0.54%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "4.374"
$ws.Range("E17").Value = "1.84%"
$ws.Range("B18").Value = "BitpandaEcosystemToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D18").Value = "0.3313"
$ws.Range("E18").Value = "0.73%"
$ws.Range("B19").Value = "MCDex"
$ws.Range("C19").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D19").Value = "8.102"
$ws.Range("E19").Value = "1.57%"
$ws.Range("B20").Value = "ProBitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D20").Value = "0.1360"
$ws.Range("E20").Value = "-2.84%"
$ws.Range("B21").Value = "ZBToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D21").Value = "0.3102"
$ws.Range("E21").Value = "11.08%"
$ws.Range("D22").Value = "0.04084"
$ws.Range("E22").Value = "3.91%"
$ws.Range("D23").Value = "0.001265"
$ws.Range("E23").Value = "0.30%"
$ws.Range("D24").Value = "0.003892"
$ws.Range("E24").Value = "2.90%"
$ws.Range("E25").Value = "-3.97%"
$ws.Range("D38").Value = "0.02417"
$ws.Range("E38").Value = "6.18%"
$ws.Range("D39").Value = "0.05204"
$ws.Range("E39").Value = "4.63%"
$ws.Range("D40").Value = "0.006301"
$ws.Range("E40").Value = "11.85%"
$ws.Range("D41").Value = "0.007779"
$ws.Range("E41").Value = "1.05%"
$ws.Range("E42").Value = "3.19%"
$ws.Range("D43").Value = "0.007392"
$ws.Range("E43").Value = "0.79%"
$ws.Range("D44").Value = "0.007256"
$ws.Range("E44").Value = "-12.48%"
$ws.Range("D45").Value = "0.2972"
$ws.Range("E45").Value = "1.96%"
$ws.Range("D46").Value = "0.00006585"
$ws.Range("E46").Value = "4.30%"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").Value = "-0.13%"
$ws.Range("D48").Value = "0.03167"
$ws.Range("E48").Value = "23.77%"
$ws.Range("D49").Value = "0.004201"
$ws.Range("E49").Value = "-0.02%"
$ws.Range("D50").Value = "0.00002100"
$ws.Range("E50").Value = "-0.13%"
$ws.Range("D51").Value = "0.0002000"
$ws.Range("E51").Value = "-0.13%"
